$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.493952035903931
$ws.Range("B1").Value = 1.496298551559448
$ws.Range("C1").Value = 1.60165548324585
$ws.Range("D1").Value = 2.276792287826538
$ws.Range("E1").Value = 4.348421573638916
